$d = $word.ActiveDocument

$replacements = @(
    @{old="615÷7="; new="646÷4="},
    @{old="908÷8="; new="525÷8="},
    @{old="214÷8="; new="195÷6="},
    @{old="732÷9="; new="570÷7="},
    @{old="874÷7="; new="747÷3="},
    @{old="746÷8="; new="613÷2="},
    @{old="333÷9="; new="394÷2="},
    @{old="694÷9="; new="704÷9="},
    @{old="862÷2="; new="233÷2="},
    @{old="194÷8="; new="413÷2="},
    @{old="549÷3="; new="238÷6="},
    @{old="995÷5="; new="993÷8="},
    @{old="197÷7="; new="217÷2="},
    @{old="431÷2="; new="859÷6="},
    @{old="841÷8="; new="436÷5="},
    @{old="364÷8="; new="898÷4="},
    @{old="334÷5="; new="627÷9="},
    @{old="798÷5="; new="989÷8="},
    @{old="559÷3="; new="117÷2="},
    @{old="814÷3="; new="109÷4="},
    @{old="881÷3="; new="105÷9="},
    @{old="473÷8="; new="432÷4="},
    @{old="833÷9="; new="919÷9="},
    @{old="570÷4="; new="847÷8="},
    @{old="932÷3="; new="758÷6="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
